$d = $word.ActiveDocument

# The "dateModified, timesSold" field list item loses its trailing
# ", timesSold" (comma + space + word), leaving just "...dateModified".
# Locate "dateModified" and "timesSold" dynamically (rather than hard-coding
# character offsets) so the edit is resilient to any earlier content.

$afterDateModified = $d.Content
$afterDateModified.Find.Execute("dateModified") | Out-Null
$startDelete = $afterDateModified.End

$afterTimesSold = $d.Content
$afterTimesSold.Find.Execute("timesSold") | Out-Null
$endDelete = $afterTimesSold.End

$d.Range($startDelete, $endDelete).Delete()
